$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update StructureDefinition summary properties ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now set to "Alvearie Team"
$meta.Range("B9").Value = "Alvearie Team"

# The old row 10 was a duplicate "Contact" / "No display for ContactDetail" row
# (same content duplicated into row 11). Turn row 10 into the new
# "Jurisdiction" / "United States of America" row, then delete the now
# redundant duplicate row 11 so the rest of the table shifts up by one.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements": update root Extension row's Short/Definition text ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Employee Average Scheduled Hours"
$elements.Range("L2").Value = "Average number of hours the employee is scheduled to work per day"
